$wb = $excel.ActiveWorkbook

# --- Update conversion note on sheet "Hoja1" (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$old1 = "1000 Bs = 9.07 = 37895.08 pesos"
$new1 = "1000 Bs = 8.93 = 37200.08 pesos"
$old2 = "37895.08 pesos = 9.06 = 974.7 Bs"
$new2 = "37200.08 pesos = 8.87 = 948.41 Bs"

$text = $ws1.Range("A1").Text
$text = $text.Replace($old1, $new1)
$text = $text.Replace($old2, $new2)
$ws1.Range("A1").Value = $text

# --- Update exchange rate figures on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 111.989
$ws2.Range("O10").Value = 4166
$ws2.Range("N12").Value = 4193
$ws2.Range("O12").Value = 106.9
